# Update the "想去人数" (interest count) values in column F for the
# "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 472
$ws1.Range("F5").Value  = 76
$ws1.Range("F6").Value  = 321
$ws1.Range("F8").Value  = 514
$ws1.Range("F11").Value = 179
$ws1.Range("F13").Value = 26
$ws1.Range("F19").Value = 617
$ws1.Range("F22").Value = 2423
$ws1.Range("F23").Value = 13
$ws1.Range("F30").Value = 2830
$ws1.Range("F33").Value = 119
$ws1.Range("F34").Value = 677
$ws1.Range("F36").Value = 1853
$ws1.Range("F38").Value = 1866

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 472
$ws4.Range("F5").Value  = 76
$ws4.Range("F8").Value  = 321
$ws4.Range("F10").Value = 514
$ws4.Range("F13").Value = 179
$ws4.Range("F15").Value = 26
$ws4.Range("F21").Value = 617
$ws4.Range("F24").Value = 2423
$ws4.Range("F28").Value = 2830
$ws4.Range("F31").Value = 119
$ws4.Range("F34").Value = 677
$ws4.Range("F36").Value = 1853
$ws4.Range("F39").Value = 1866

$wb.Save()
